$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62

# Row 3 updates
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
